# Collapse the word-by-word text runs (e.g. "Blank" " " "slides" ...)
# back down into a single run per paragraph, as described by the diff.

$p = $ppt.ActivePresentation

function Set-ShapeText($shape, $text) {
    # The host treats an assignment that produces the same *concatenated*
    # string as a no-op and leaves the original (split) runs untouched, so
    # first push a throwaway value to force a genuine content change, then
    # set the real text. That makes PowerPoint rebuild the paragraph as a
    # single run.
    $shape.TextFrame.TextRange.Text = "."
    $shape.TextFrame.TextRange.Text = $text
}

# Slide 1: "Section Header (with background image)"
$s1 = $p.Slides.Item(1)
Set-ShapeText $s1.Shapes.Item("Title 1") "Section Header (with background image)"

# Slide 2: "Slide 1"
$s2 = $p.Slides.Item(2)
Set-ShapeText $s2.Shapes.Item("Title 1") "Slide 1"

# Slide 3: "Slide 2"
$s3 = $p.Slides.Item(3)
Set-ShapeText $s3.Shapes.Item("Title 1") "Slide 2"

# Slide 4: "Slide 3"
$s4 = $p.Slides.Item(4)
Set-ShapeText $s4.Shapes.Item("Title 1") "Slide 3"

# Slide 5: "Slide 4" title and "An image" caption textbox
$s5 = $p.Slides.Item(5)
Set-ShapeText $s5.Shapes.Item("Title 1") "Slide 4"
Set-ShapeText $s5.Shapes.Item("TextBox 3") "An image"

# Slide 6's notes page: "Blank slides can have background images."
$s6 = $p.Slides.Item(6)
$notes = $s6.NotesPage
Set-ShapeText $notes.Shapes.Item("Notes Placeholder 2") "Blank slides can have background images."
